$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.837.20'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.874.34'
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '527.22'
$ws.Range('E5').Value = '  +7.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.84'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.607'
$ws.Range('E7').Value = '  -2.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.714'
$ws.Range('E9').Value = '  -3.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.169'
$ws.Range('E10').Value = '  -5.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000327'
$ws.Range('E11').Value = '  -6.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '41.90'
$ws.Range('E12').Value = '  -2.60%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.25'
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.482.13'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.64'
$ws.Range('E15').Value = '  +8.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.873.24'
$ws.Range('E16').Value = '  -1.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.21'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('E18').Value = '  +6.20%  '
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.806.72'
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '422.58'
$ws.Range('E21').Value = '  -2.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.41'
$ws.Range('E22').Value = '  -2.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.08'
$ws.Range('E23').Value = '  -4.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.27'
$ws.Range('E24').Value = '  -3.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.00'
$ws.Range('E25').Value = '  +6.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.33'
$ws.Range('E26').Value = '  -7.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.66'
$ws.Range('E27').Value = '  -3.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.08'
$ws.Range('E28').Value = '  -3.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '687.47'
$ws.Range('E29').Value = '  -3.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '13.13'
$ws.Range('E30').Value = '  -1.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.126'
$ws.Range('E31').Value = '  -3.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.84'
$ws.Range('E32').Value = '  -2.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '66.97'
$ws.Range('E33').Value = '  +8.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.436'
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.92'
$ws.Range('E35').Value = '  -4.19%  '
$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0850'
$ws.Range('E36').Value = '  -3.42%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '39.95'
$ws.Range('E37').Value = '  -2.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.22'
$ws.Range('E41').Value = '  +2.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.20'
$ws.Range('E42').Value = '  +5.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0479'
$ws.Range('E43').Value = '  -2.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.75'
$ws.Range('E44').Value = '  -9.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.39'
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.805.51'
$ws.Range('E47').Value = '  +15.79%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.95'
$ws.Range('E48').Value = '  +4.80%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000270'
$ws.Range('E49').Value = '  +13.32%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₆0343'
$ws.Range('E50').Value = '  -8.85%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.53'
$ws.Range('E51').Value = '  +1.32%  '
